$d = $word.ActiveDocument

# --- Locate the first empty paragraph right after "OR-Mapping " -----
# (rather than hard-coding a paragraph index, find it by content so
# the script is resilient to minor shifts elsewhere in the document).
$anchor = $d.Content
$anchor.Find.Execute("OR-Mapping", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0) | Out-Null
$anchorEnd = $anchor.End

$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -ge $anchorEnd) {
        $targetIndex = $i
        break
    }
}

# --- 1) Move the _GoBack bookmark off of its current paragraph ------
# It currently sits on the paragraph that, after "Speedikon DAMS...",
# just holds a single space. _GoBack is reachable directly by name
# even though Word normally hides it from the Bookmarks collection.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# --- 2) Give that first empty paragraph the new "Zxing" text --------
$p = $d.Paragraphs.Item($targetIndex)
$p.Range.Text = "ZxingX"

# --- 3) Re-create _GoBack right after "Zxing", collapsed ------------
# (mirrors bookmarkStart immediately followed by bookmarkEnd, i.e. a
# zero-length bookmark). Bookmarks.Add on a genuinely collapsed range
# isn't reliable in this runtime, so temporarily wrap a trailing
# placeholder character with the bookmark, then delete that character;
# the bookmark collapses down to sit right after "Zxing".
$p = $d.Paragraphs.Item($targetIndex)
$placeholder = $d.Range($p.Range.End - 2, $p.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $placeholder) | Out-Null
$placeholder2 = $d.Range($p.Range.End - 2, $p.Range.End - 1)
$placeholder2.Text = ""
